$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text tweaks: bump the report Volume/Number and the week-covering
#    dates. These live as rich-text runs inside a shared string, so we use
#    Characters() to replace just the affected substring in place.
# ---------------------------------------------------------------------------

function Replace-Substring($range, [string]$old, [string]$new) {
    $full = [string]$range.Value2
    $idx = $full.IndexOf($old)
    if ($idx -ge 0) {
        $chars = $range.Characters($idx + 1, $old.Length)
        $chars.Text = $new
    }
}

Replace-Substring $ws.Range("A8") "51" "52"
Replace-Substring $ws.Range("C9") "12/18/2023" "12/25/2023"
Replace-Substring $ws.Range("C9") "12/24/2023" "12/31/2023"

# ---------------------------------------------------------------------------
# 2) Helper to change a cell that currently holds the text placeholder "0"
#    (shared-string style, s="14") into a real number while reusing the
#    numeric-cell style (s="15") - mirrored by copying the format from a
#    known numeric donor cell so no new style entries are created.
# ---------------------------------------------------------------------------
function Set-AsNumber($addr, $value) {
    $donor = $ws.Range("G15")
    $dst = $ws.Range($addr)
    $donor.Copy($dst)
    $dst.Value2 = $value
}

# ...and the inverse: change a numeric cell (s="15") back into the shared
# text placeholder "0" (s="14"), by copying format+value from a donor cell
# that already stores that placeholder.
function Set-AsZeroPlaceholder($addr) {
    $donor = $ws.Range("D14")
    $dst = $ws.Range($addr)
    $donor.Copy($dst)
}

# ---------------------------------------------------------------------------
# 3) Precinct crime-stat table refresh (rows 15-29) - new weekly figures.
# ---------------------------------------------------------------------------

# Row 15 - Rape
Set-AsNumber "C15" 2
$ws.Range("F15").Value2 = 5
$ws.Range("H15").Value2 = 25
$ws.Range("I15").Value2 = 38
$ws.Range("K15").Value2 = -5
$ws.Range("L15").Value2 = -7.317073170731
$ws.Range("M15").Value2 = 35.714285714285
$ws.Range("N15").Value2 = 11.764705882352

# Row 16 - Robbery
$ws.Range("C16").Value2 = 21
$ws.Range("E16").Value2 = 425
$ws.Range("F16").Value2 = 52
$ws.Range("G16").Value2 = 18
$ws.Range("H16").Value2 = 188.888888888889
$ws.Range("I16").Value2 = 479
$ws.Range("J16").Value2 = 399
$ws.Range("K16").Value2 = 20.050125313283
$ws.Range("L16").Value2 = 85.658914728682
$ws.Range("M16").Value2 = 32.320441988950
$ws.Range("N16").Value2 = -68.916288124594

# Row 17 - Fel. Assault
$ws.Range("C17").Value2 = 10
$ws.Range("D17").Value2 = 6
$ws.Range("E17").Value2 = 66.666666666666
$ws.Range("F17").Value2 = 44
$ws.Range("G17").Value2 = 28
$ws.Range("H17").Value2 = 57.142857142857
$ws.Range("I17").Value2 = 732
$ws.Range("J17").Value2 = 525
$ws.Range("K17").Value2 = 39.428571428571
$ws.Range("L17").Value2 = 72.235294117647
$ws.Range("M17").Value2 = 187.058823529412
$ws.Range("N17").Value2 = 47.283702213279

# Row 18 - Burglary
$ws.Range("C18").Value2 = 3
$ws.Range("D18").Value2 = 9
$ws.Range("E18").Value2 = -66.666666666666
$ws.Range("F18").Value2 = 9
$ws.Range("G18").Value2 = 11
$ws.Range("H18").Value2 = -18.181818181818
$ws.Range("I18").Value2 = 188
$ws.Range("J18").Value2 = 163
$ws.Range("K18").Value2 = 15.337423312883
$ws.Range("L18").Value2 = 1.621621621621
$ws.Range("M18").Value2 = -35.836177474402
$ws.Range("N18").Value2 = -91.316397228637

# Row 19 - Gr. Larceny
$ws.Range("C19").Value2 = 24
$ws.Range("D19").Value2 = 20
$ws.Range("E19").Value2 = 20
$ws.Range("F19").Value2 = 84
$ws.Range("G19").Value2 = 67
$ws.Range("H19").Value2 = 25.373134328358
$ws.Range("I19").Value2 = 1075
$ws.Range("J19").Value2 = 1155
$ws.Range("K19").Value2 = -6.926406926406
$ws.Range("L19").Value2 = 5.911330049261
$ws.Range("M19").Value2 = 92.998204667863
$ws.Range("N19").Value2 = -16.015625

# Row 20 - G.L.A.
$ws.Range("C20").Value2 = 4
$ws.Range("D20").Value2 = 2
$ws.Range("E20").Value2 = 100
$ws.Range("F20").Value2 = 22
$ws.Range("G20").Value2 = 18
$ws.Range("H20").Value2 = 22.222222222222
$ws.Range("I20").Value2 = 336
$ws.Range("J20").Value2 = 234
$ws.Range("K20").Value2 = 43.589743589743
$ws.Range("L20").Value2 = 80.645161290322
$ws.Range("M20").Value2 = 133.333333333333
$ws.Range("N20").Value2 = -85.403996524761

# Row 21 - TOTAL
$ws.Range("C21").Value2 = 64
$ws.Range("D21").Value2 = 41
$ws.Range("E21").Value2 = 56.097560975609
$ws.Range("F21").Value2 = 216
$ws.Range("G21").Value2 = 146
$ws.Range("H21").Value2 = 47.945205479452
$ws.Range("I21").Value2 = 2851
$ws.Range("J21").Value2 = 2524
$ws.Range("K21").Value2 = 12.955625990491
$ws.Range("L21").Value2 = 34.862819299905
$ws.Range("M21").Value2 = 73.524041387705
$ws.Range("N21").Value2 = -63.672273190621

# Row 22 - Transit
Set-AsNumber "C22" 1
$ws.Range("E22").Value2 = 0
$ws.Range("F22").Value2 = 3
$ws.Range("H22").Value2 = 50
$ws.Range("I22").Value2 = 61
$ws.Range("J22").Value2 = 44
$ws.Range("K22").Value2 = 38.636363636363
$ws.Range("L22").Value2 = 117.857142857143
$ws.Range("M22").Value2 = 96.774193548387

# Row 24 - Petit Larceny
$ws.Range("C24").Value2 = 69
$ws.Range("D24").Value2 = 53
$ws.Range("E24").Value2 = 30.188679245283
$ws.Range("F24").Value2 = 263
$ws.Range("G24").Value2 = 236
$ws.Range("H24").Value2 = 11.440677966101
$ws.Range("I24").Value2 = 2875
$ws.Range("J24").Value2 = 2453
$ws.Range("K24").Value2 = 17.203424378312
$ws.Range("L24").Value2 = 56.25
$ws.Range("M24").Value2 = 57.189721159103

# Row 25 - Misd. Assault
$ws.Range("C25").Value2 = 32
$ws.Range("D25").Value2 = 21
$ws.Range("E25").Value2 = 52.380952380952
$ws.Range("F25").Value2 = 91
$ws.Range("G25").Value2 = 68
$ws.Range("H25").Value2 = 33.823529411764
$ws.Range("I25").Value2 = 1179
$ws.Range("J25").Value2 = 912
$ws.Range("K25").Value2 = 29.276315789473
$ws.Range("L25").Value2 = 46.459627329192
$ws.Range("M25").Value2 = 83.074534161490

# Row 26 - UCR Rape*
Set-AsNumber "C26" 2
$ws.Range("I26").Value2 = 61
$ws.Range("K26").Value2 = 3.389830508474
$ws.Range("L26").Value2 = -4.6875

# Row 27 - Other Sex Crimes
$ws.Range("C27").Value2 = 2
$ws.Range("D27").Value2 = 2
$ws.Range("E27").Value2 = 0
$ws.Range("F27").Value2 = 9
$ws.Range("G27").Value2 = 6
$ws.Range("H27").Value2 = 50
$ws.Range("I27").Value2 = 132
$ws.Range("J27").Value2 = 127
$ws.Range("K27").Value2 = 3.937007874015
$ws.Range("L27").Value2 = 37.5

# Row 28 - Shooting Vic.
Set-AsZeroPlaceholder "F28"
$ws.Range("H28").Value2 = -100
$ws.Range("N28").Value2 = -80.645161290322

# Row 29 - Shooting Inc.
Set-AsZeroPlaceholder "F29"
$ws.Range("H29").Value2 = -100
$ws.Range("N29").Value2 = -81.818181818181
